$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete row 2 (A2=45292, B2=7858) and shift remaining rows up,
# keeping only five-year intervals starting at 45839.
$ws.Rows.Item(2).Delete()
